$wb = $excel.ActiveWorkbook

# Update the "zh-cn" sheet's first data row (row 2) with newly generated
# handoff/handback report timestamps.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 10:50:07"
$wsZhCn.Range("H2").Value = "2016-03-22 10:50:30"

# Update the "de-de" sheet's first data row (row 2) with newly generated
# handoff/handback report timestamps.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 10:50:11"
$wsDeDe.Range("H2").Value = "2016-03-22 10:50:40"
